$wb = $excel.ActiveWorkbook

# --- Yearly sheet: update 2017 "Apr" taxable dividend (row 6, column L) ---
# O6, L15 and O15 are formulas (SUM) and recalc automatically.
$wsYearly = $wb.Worksheets.Item("Yearly")
$wsYearly.Range("L6").Value = 44.12
$wsYearly.Range("O6").Select()

# --- All Time sheet: update view selection / scroll position ---
# F8/I8 (=Yearly!L15 / shared SUM) and F46/I46 (=SUM(F6:F45) / shared SUM)
# recalc automatically from the Yearly-sheet edit above.
$wsAllTime = $wb.Worksheets.Item("All Time")
$wsAllTime.Activate()
$wsAllTime.Application.ActiveWindow.ScrollRow = 31
$wsAllTime.Application.ActiveWindow.ScrollColumn = 1
$wsAllTime.Range("K39").Select()
